# act 08-08-2025 se agregan mas equipos y nueva funcionalidad
#
# Fills in the previously-empty match row 23 (Melgar vs Juan Pablo II,
# 05/08/2025) with the new game's stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A23 is a date written as literal text (matching the rest of column A,
# which stores "dd/mm/yyyy" as plain strings rather than real dates).
# Prefix with an apostrophe so Excel doesn't auto-convert it to a date
# serial, then copy the neighboring cell's (default) style back over so
# the quote-prefix flag doesn't leave a stray style behind.
$ws.Cells.Item(23, 1).Value = "'05/08/2025"
$ws.Cells.Item(23, 1).Style = $ws.Cells.Item(22, 1).Style

$ws.Cells.Item(23, 2).Value = "Melgar"
$ws.Cells.Item(23, 3).Value = 1
$ws.Cells.Item(23, 4).Value = 1
$ws.Cells.Item(23, 5).Value = "Juan Pablo II"
$ws.Cells.Item(23, 6).Value = "D"
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 1
$ws.Cells.Item(23, 9).Value = 1
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 0.17
$ws.Cells.Item(23, 12).Value = 1.38
$ws.Cells.Item(23, 13).Value = 5
$ws.Cells.Item(23, 14).Value = 21
$ws.Cells.Item(23, 15).Value = 2
$ws.Cells.Item(23, 16).Value = 5
